$wb = $excel.ActiveWorkbook

$changes = @{
    "InversTeilsortiert1000.dat" = @{
        "B2" = 7478600.0
        "B3" = 7504200.0
        "C3" = 32104.0
        "B4" = 1295500.0
        "B5" = 68500.0
        "C5" = 32128.0
        "B6" = 2801600.0
        "C6" = 32128.0
        "B7" = 1733500.0
        "D7" = 23262.0
        "E7" = 5754.0
        "B8" = 4301100.0
    }
    "InversTeilsortiert10000.dat" = @{
        "B2" = 178647100.0
        "B3" = 189554900.0
        "C3" = 320104.0
        "B4" = 3896100.0
        "B5" = 1034700.0
        "C5" = 320128.0
        "B6" = 173730700.0
        "C6" = 320128.0
        "B7" = 3305900.0
        "D7" = 325920.0
        "E7" = 86640.0
        "B8" = 113204300.0
    }
    "InversTeilsortiert100000.dat" = @{
        "B2" = -423809488.0
        "B3" = 1112973612.0
        "C3" = 3200104.0
        "B4" = 11716600.0
        "B5" = 3169000.0
        "C5" = 3200128.0
        "B6" = -1426433696.0
        "C6" = 3200128.0
        "B7" = 12753700.0
        "D7" = 4000593.0
        "E7" = 1111531.0
        "B8" = 1349283904.0
    }
    "Random1000.dat" = @{
        "B2" = 1513600.0
        "B3" = 2033900.0
        "C3" = 32104.0
        "B4" = 165600.0
        "C5" = 32128.0
        "B6" = 348500.0
        "C6" = 32128.0
        "B7" = 225100.0
        "D7" = 4020549.0
        "E7" = 1116183.0
        "B8" = 1053500.0
    }
    "Random10000.dat" = @{
        "B2" = 157672200.0
        "B3" = 157194900.0
        "C3" = 320104.0
        "B4" = 1907500.0
        "B5" = 109100.0
        "C5" = 320128.0
        "B6" = 34066100.0
        "C6" = 320128.0
        "B7" = 2505900.0
        "D7" = 4317285.0
        "E7" = 1195095.0
        "B8" = 77885900.0
    }
    "Random100000.dat" = @{
        "B2" = -1400865480.0
        "B3" = 781544112.0
        "C3" = 3200104.0
        "B4" = 13172600.0
        "B5" = 241500.0
        "C5" = 3200128.0
        "B6" = 2037912500.0
        "C6" = 3200128.0
        "B7" = 17827200.0
        "D7" = 8014272.0
        "E7" = 2227424.0
        "B8" = 1528776804.0
    }
    "Teilsortiert1000.dat" = @{
        "B2" = 583400.0
        "B3" = 283100.0
        "C3" = 32104.0
        "B4" = 102600.0
        "B5" = 1900.0
        "C5" = 32128.0
        "B6" = 3349800.0
        "C6" = 32128.0
        "B7" = 153300.0
        "D7" = 8035692.0
        "E7" = 2232564.0
        "B8" = 580900.0
    }
    "Teilsortiert10000.dat" = @{
        "B2" = 49833800.0
        "B3" = 34998800.0
        "C3" = 320104.0
        "B4" = 1397300.0
        "B5" = 46300.0
        "C5" = 320128.0
        "B6" = 59009300.0
        "C6" = 320128.0
        "B7" = 2109400.0
        "D7" = 8328342.0
        "E7" = 2310114.0
        "B8" = 67033100.0
    }
    "Teilsortiert100000.dat" = @{
        "B2" = -345289692.0
        "B3" = -749845396.0
        "C3" = 3200104.0
        "B4" = 15407200.0
        "B5" = 215100.0
        "C5" = 3200128.0
        "B6" = -714967496.0
        "C6" = 3200128.0
        "B7" = 12240500.0
        "D7" = 11741022.0
        "E7" = 3247674.0
        "B8" = 1557466204.0
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($changes.ContainsKey($name)) {
        $cellChanges = $changes[$name]
        foreach ($addr in $cellChanges.Keys) {
            $ws.Range($addr).Value = $cellChanges[$addr]
        }
    }
}
